# Updates cryptos list values (Price / Volume(1h) columns, plus a couple
# of re-ranked coin rows) to match the latest scrape, per commit message:
# "Updated cryptos list on Mon Mar 20 16:16:00 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.859.14'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '1.775.00'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('D4').Value = '''1.008'
$ws.Range('E4').Value = '  +0.70%  '
$ws.Range('D5').Value = '''339.84'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = '''1.005'
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('D7').Value = '''0.3825'
$ws.Range('E7').Value = '  -2.60%  '
$ws.Range('D8').Value = '''0.3412'
$ws.Range('E8').Value = '  -1.34%  '
$ws.Range('D9').Value = '''46.86'
$ws.Range('E9').Value = '  -2.45%  '
$ws.Range('D10').Value = '''1.140'
$ws.Range('E10').Value = '  -4.91%  '
$ws.Range('D11').Value = '''0.07405'
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('D12').Value = '''23.32'
$ws.Range('E12').Value = '  +5.46%  '
$ws.Range('D13').Value = '''1.009'
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('D14').Value = '''6.376'
$ws.Range('E14').Value = '  -2.13%  '
$ws.Range('D15').Value = '''7.413'
$ws.Range('E15').Value = '  +3.54%  '
$ws.Range('D16').Value = '1.778.54'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').Value = '''0.00001076'
$ws.Range('E17').Value = '  -2.27%  '
$ws.Range('D18').Value = '''0.06709'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('E19').Value = '  -2.91%  '
$ws.Range('D20').Value = '''1.003'
$ws.Range('E20').Value = '  +0.50%  '
$ws.Range('D21').Value = '''17.41'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').Value = '''6.416'
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('D23').Value = '27.859.99'
$ws.Range('E23').Value = '  +0.97%  '
$ws.Range('D24').Value = '''12.04'
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('D25').Value = '''2.406'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').Value = '''1.453'
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('E27').Value = '  -2.36%  '
$ws.Range('D28').Value = '''2.417'
$ws.Range('E28').Value = '  -3.79%  '
$ws.Range('D29').Value = '''154.21'
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('D30').Value = '1.982.01'
$ws.Range('E30').Value = '  -0.97%  '
$ws.Range('D31').Value = '''134.02'
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('D32').Value = '''4.034'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('D33').Value = '''6.038'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').Value = '''0.08849'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('D35').Value = '''12.68'
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '''0.02393'
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').Value = '''0.6817'
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('D38').Value = '''0.06389'
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('D39').Value = '''5.290'
$ws.Range('E39').Value = '  -3.10%  '
$ws.Range('D40').Value = '''0.2155'
$ws.Range('E40').Value = '  -2.52%  '
$ws.Range('B41').Value = 'WEMIXTOKEN'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = '''1.502'
$ws.Range('E41').Value = '  -6.85%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''1.234'
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').Value = '''8.187'
$ws.Range('E43').Value = '  -2.85%  '
$ws.Range('D44').Value = '''14.19'
$ws.Range('E44').Value = '  -3.63%  '
$ws.Range('D45').Value = '''1.004'
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').Value = '''0.6223'
$ws.Range('E46').Value = '  -3.40%  '
$ws.Range('D47').Value = '''3.868'
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('D48').Value = '''134.09'
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('D49').Value = '''2.065'
$ws.Range('E49').Value = '  -3.87%  '
$ws.Range('D50').Value = '''0.07420'
$ws.Range('E50').Value = '  +2.96%  '
$ws.Range('D51').Value = '''1.205'
$ws.Range('E51').Value = '  +4.14%  '
